# Applies two edits to the design doc:
#  1. The empty paragraph right after the "she codes" logo image gets a
#     new bold run "Sapir mizrahi" (and the paragraph mark itself picks
#     up the same run formatting: bold, size 12pt, en-US).
#  2. The mis-spelled / flagged variable name "ne_max" is split into two
#     runs "max" + "_ne" (i.e. the visible text becomes "max_ne").
#
# Both edits are applied via Range.InsertXML() over the *entire* target
# paragraph (Paragraph.Range, which includes the end-of-paragraph mark),
# supplying a complete <w:p> element - including the paragraph's original
# identity attributes (w14:paraId / rsid*) - so InsertXML performs a clean
# whole-paragraph swap instead of merely appending runs to the paragraph
# (which is what happens if the replacement XML is not itself a <w:p>, or
# if the target Range does not span the full paragraph).

$d = $word.ActiveDocument
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'
$w14Ns = 'xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"'
$count = $d.Paragraphs.Count

# ---------------------------------------------------------------------
# Edit 1: add "Sapir mizrahi" to the empty paragraph under the logo.
# Locate it as "the paragraph right after the one holding the inline
# picture" rather than a hard-coded index.
# ---------------------------------------------------------------------
$nameIndex = -1
for ($i = 1; $i -le $count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.InlineShapes.Count -gt 0) {
        $nameIndex = $i + 1
        break
    }
}
if ($nameIndex -eq -1) {
    $nameIndex = 2
}
$namePara = $d.Paragraphs.Item($nameIndex)

$nameXml = '<w:p ' + $wNs + ' ' + $w14Ns + ' w14:paraId="48550C09" w14:textId="77777777" w:rsidR="0046208C" w:rsidRDefault="0046208C" w:rsidP="001E7EA5">' + `
    '<w:pPr>' + `
        '<w:spacing w:line="360" w:lineRule="auto"/>' + `
        '<w:rPr>' + `
            '<w:b/>' + `
            '<w:bCs/>' + `
            '<w:sz w:val="24"/>' + `
            '<w:szCs w:val="24"/>' + `
            '<w:lang w:val="en-US"/>' + `
        '</w:rPr>' + `
    '</w:pPr>' + `
    '<w:r>' + `
        '<w:rPr>' + `
            '<w:b/>' + `
            '<w:bCs/>' + `
            '<w:sz w:val="24"/>' + `
            '<w:szCs w:val="24"/>' + `
            '<w:lang w:val="en-US"/>' + `
        '</w:rPr>' + `
        '<w:t>Sapir mizrahi</w:t>' + `
    '</w:r>' + `
'</w:p>'

$namePara.Range.InsertXML($nameXml)

# ---------------------------------------------------------------------
# Edit 2: split the run "ne_max" into "max" + "_ne".
# ---------------------------------------------------------------------
$neMaxIndex = -1
for ($i = 1; $i -le $count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text -like "*ne_max*") {
        $neMaxIndex = $i
        break
    }
}
$neMaxPara = $d.Paragraphs.Item($neMaxIndex)

$neMaxXml = '<w:p ' + $wNs + ' ' + $w14Ns + ' w14:paraId="02AEF142" w14:textId="77777777" w:rsidR="001E7EA5" w:rsidRDefault="001E7EA5" w:rsidP="001E7EA5">' + `
    '<w:pPr>' + `
        '<w:pStyle w:val="ListParagraph"/>' + `
        '<w:numPr><w:ilvl w:val="0"/><w:numId w:val="8"/></w:numPr>' + `
        '<w:spacing w:line="360" w:lineRule="auto"/>' + `
        '<w:rPr>' + `
            '<w:sz w:val="24"/>' + `
            '<w:szCs w:val="24"/>' + `
            '<w:lang w:val="en-US"/>' + `
        '</w:rPr>' + `
    '</w:pPr>' + `
    '<w:r>' + `
        '<w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="en-US"/></w:rPr>' + `
        '<w:t xml:space="preserve">The most efficient number of trees (variable </w:t>' + `
    '</w:r>' + `
    '<w:proofErr w:type="spellStart"/>' + `
    '<w:r>' + `
        '<w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="en-US"/></w:rPr>' + `
        '<w:t>max</w:t>' + `
    '</w:r>' + `
    '<w:r>' + `
        '<w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="en-US"/></w:rPr>' + `
        '<w:t>_ne</w:t>' + `
    '</w:r>' + `
    '<w:proofErr w:type="spellEnd"/>' + `
    '<w:r>' + `
        '<w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="en-US"/></w:rPr>' + `
        '<w:t>) was used in building the model.</w:t>' + `
    '</w:r>' + `
'</w:p>'

$neMaxPara.Range.InsertXML($neMaxXml)

Write-Output "Applied both edits."
